$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 355 ---
$ws.Range("D355").Value = 50
$ws.Range("E355").Value = 52.5

# --- Row 356 ---
$ws.Range("D356").Value = 98
$ws.Range("E356").Value = 45
$ws.Range("H356").Value = 16.1

# --- Row 357 (adds D357) ---
$ws.Range("D357").Value = 37
$ws.Range("E357").Value = 42.6

# --- Row 358 (adds D358) ---
$ws.Range("D358").Value = 41
$ws.Range("E358").Value = 46.4

# --- Row 359 ---
$ws.Range("E359").Value = 59.4
$ws.Range("F359").Value = 0.5
$ws.Range("H359").Value = 13.2

# --- Row 360 ---
$ws.Range("E360").Value = 57.9
$ws.Range("F360").Value = 0.5
$ws.Range("H360").Value = 18.2

# --- Row 361 ---
$ws.Range("E361").Value = 39.9
$ws.Range("F361").Value = 0
$ws.Range("H361").Value = 12.8

# --- Row 362 ---
$ws.Range("E362").Value = 43.7
$ws.Range("H362").Value = 10.2

# --- Row 363 ---
$ws.Range("E363").Value = 52.5
$ws.Range("H363").Value = 12.6

# --- Row 364 ---
$ws.Range("E364").Value = 50.9
$ws.Range("F364").Value = 0.2
$ws.Range("H364").Value = 12.9

# --- Row 365 ---
$ws.Range("E365").Value = 42.4
$ws.Range("F365").Value = 0.5
$ws.Range("G365").Value = 1
$ws.Range("H365").Value = 13.4

# --- Row 366 (new) ---
$ws.Range("A366").Value = 6
$ws.Range("B366").Value = 365
$ws.Range("C366").Value = 45737
$ws.Range("E366").Value = 40.5
$ws.Range("F366").Value = 0
$ws.Range("G366").Value = 0
$ws.Range("H366").Value = 11.9

# --- Row 367 (new) ---
$ws.Range("A367").Value = 7
$ws.Range("B367").Value = 366
$ws.Range("C367").Value = 45738
$ws.Range("E367").Value = 40.8
$ws.Range("F367").Value = 0
$ws.Range("G367").Value = 0
$ws.Range("H367").Value = 11

# Copy the date-number-format from an existing date cell onto the two new
# date cells (C366/C367) so they reuse the same style index instead of
# Excel materialising a brand-new custom numFmt.
$ws.Range("C354").Copy()
$ws.Range("C366:C367").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View / selection update ---
$ws.Range("G364").Select()
